# Adds new testscript rows WAT44 / WAT45 to the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 25 carries the same border/wrap formatting (no extra row height) that the
# two new rows end up with, so re-use it as the formatting template for the
# freshly appended rows 32 and 33.
$ws.Range("A25:E25").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Range("A33:E33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New test case WAT44
$ws.Range("A32").Value = "WAT44"
$ws.Range("A33").Value = "WAT45"
$ws.Range("B32").Value = "WAT-203"
$ws.Range("C32").Value = "Verify the consistent presence of LOGO at the top left of page."
$ws.Range("C33").Value = "Verify that clicking on the LOGO navigates the user to http://clarivate.com/"
$ws.Range("B33").Value = "WAT-211"
$ws.Range("D32").Value = "Y"
$ws.Range("D33").Value = "Y"

# Restore the view state (selection) that the author ended up with.
$ws.Activate()
$null = $ws.Range("L43").Select()
